# The dataset was regenerated: the first 7 original rows (old rows 2-8)
# were dropped, the remaining 7 original rows (old rows 16-22) shifted up
# to become the new rows 2-8, and 13 freshly generated rows were appended
# as new rows 9-21. Net effect: sheet shrinks from A1:C22 to A1:C21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data block (row 1 header stays untouched).
$ws.Range("A2:C22").Clear() | Out-Null

# Final dataset (rows 2-21).
$data = @(
    @(2.08137059211731, -2.974608421325684, 1.146135926246643),
    @(-0.7811439037322998, 1.68369734287262, -1.490358471870422),
    @(-0.4285219609737396, 1.809230089187622, -1.255480766296387),
    @(-1.154382586479187, 3.296534299850464, -1.97553825378418),
    @(-1.183093309402466, 0.8458956480026245, -1.687667965888977),
    @(-0.4327980279922485, 2.791501522064209, -2.889545440673828),
    @(-2.166738986968994, 1.760666370391846, -2.120465993881226),
    @(1.305724501609802, 0.0389426611363887, -1.967597007751465),
    @(2.605340242385864, -0.3645338416099548, -0.5674937963485718),
    @(-0.7612907886505127, 2.151620149612427, 0.1557706445455551),
    @(0.5158756971359253, -2.014939069747925, 1.058934926986694),
    @(-0.7177666425704956, 0.8080220222473145, 1.101084589958191),
    @(0.3738495409488678, -2.970790386199951, 2.029447078704834),
    @(1.252579212188721, -2.954449653625488, 2.593122959136963),
    @(1.917964100837708, -2.736523628234864, 1.525177836418152),
    @(-0.9870055317878724, -1.847256541252136, 0.9668469429016112),
    @(1.084591269493103, -5.03795862197876, 1.1534663438797),
    @(-1.704924941062927, -0.2585487067699432, -1.87276017665863),
    @(0.7817547917366028, 1.986839175224304, -2.802802562713623),
    @(0.7434229850769043, 0.6293439269065857, -2.307848930358887)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
